# Update "想去人数" (interested-count) figures in the 展览 and 全部类型
# sheets to match the freshly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7
$ws1.Range("F5").Value  = 13704
$ws1.Range("F7").Value  = 37
$ws1.Range("F8").Value  = 1757
$ws1.Range("F12").Value = 44
$ws1.Range("F13").Value = 12
$ws1.Range("F15").Value = 13705
$ws1.Range("F17").Value = 611
$ws1.Range("F20").Value = 8159
$ws1.Range("F31").Value = 403
$ws1.Range("F35").Value = 387
$ws1.Range("F37").Value = 4930

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7
$ws4.Range("F3").Value  = 1041
$ws4.Range("F5").Value  = 13704
$ws4.Range("F7").Value  = 37
$ws4.Range("F8").Value  = 1757
$ws4.Range("F12").Value = 44
$ws4.Range("F13").Value = 12
$ws4.Range("F15").Value = 13705
$ws4.Range("F17").Value = 611
$ws4.Range("F20").Value = 8159
$ws4.Range("F33").Value = 403
$ws4.Range("F37").Value = 387
$ws4.Range("F39").Value = 4931
